$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New experiment data for columns B:J, rows 2-91.
# (Column A sequence numbers 0-89 already correct/unchanged for these rows.)
# Cell-by-cell assignment is used because bulk 2-D Range.Value array writes
# are not reliable in this COM host.
$rows = @(
  @('EM','EM_Q100',0,'MDUAL',2159.55,16.9,80,19,25),
  @('EM','EM_Q100',0,'SOP',15157.63,642.1,684,19,25),
  @('EM','EM_Q100',0,'pMCSKY',15176.31,601,663,19,25),
  @('EM','EM_Q100',0.2,'MDUAL',2303.36,16,78,123,123),
  @('EM','EM_Q100',0.2,'SOP',40302.55,663,745,123,123),
  @('EM','EM_Q100',0.2,'pMCSKY',40450.24,606.5,687,123,123),
  @('EM','EM_Q100',0.4,'MDUAL',2557.97,16.1,79,306,307),
  @('EM','EM_Q100',0.4,'SOP',42192.1,672.4,764,306,307),
  @('EM','EM_Q100',0.4,'pMCSKY',40617.45,615.7,707,306,307),
  @('EM','EM_Q100',0.8,'MDUAL',1783.88,15.9,78,26,26),
  @('EM','EM_Q100',0.8,'SOP',44242.53,671.6,763,26,26),
  @('EM','EM_Q100',0.8,'pMCSKY',43280.23,614.8,705,26,26),
  @('EM','EM_Q100',1,'MDUAL',2297.65,15.8,78,110,112),
  @('EM','EM_Q100',1,'SOP',41559.51,666.7,756,110,112),
  @('EM','EM_Q100',1,'pMCSKY',42070.2,609.7,699,110,112),
  @('FC','FC_Q100',0,'MDUAL',3151.88,21.1,57,330,729),
  @('FC','FC_Q100',0,'SOP',4739.57,228.3,332,330,729),
  @('FC','FC_Q100',0,'pMCSKY',4631.97,187.4,289,330,729),
  @('FC','FC_Q100',0.2,'MDUAL',2039.48,20.9,74,313,615),
  @('FC','FC_Q100',0.2,'SOP',51983.87,220.3,319,313,615),
  @('FC','FC_Q100',0.2,'pMCSKY',52842.81,164.2,266,313,615),
  @('FC','FC_Q100',0.4,'MDUAL',2380.05,21,59,323,646),
  @('FC','FC_Q100',0.4,'SOP',25678.27,222.6,332,323,646),
  @('FC','FC_Q100',0.4,'pMCSKY',26204.89,166.5,274,323,646),
  @('FC','FC_Q100',0.8,'MDUAL',2428.62,21,54,333,714),
  @('FC','FC_Q100',0.8,'SOP',25087.03,222.2,328,333,714),
  @('FC','FC_Q100',0.8,'pMCSKY',25673.3,166,270,333,714),
  @('FC','FC_Q100',1,'MDUAL',2475.98,21,45,351,736),
  @('FC','FC_Q100',1,'SOP',26038.8,223.6,322,351,736),
  @('FC','FC_Q100',1,'pMCSKY',26698.3,167.5,264,351,736),
  @('GAS','GAS_Q100',0,'MDUAL',560.43,19,80,34,83),
  @('GAS','GAS_Q100',0,'SOP',9373.28,346.6,419,34,83),
  @('GAS','GAS_Q100',0,'pMCSKY',8732.98,303.5,377,34,83),
  @('GAS','GAS_Q100',0.2,'MDUAL',417.96,20.2,80,32,79),
  @('GAS','GAS_Q100',0.2,'SOP',32568.16,373.7,453,32,79),
  @('GAS','GAS_Q100',0.2,'pMCSKY',32467.03,317,396,32,79),
  @('GAS','GAS_Q100',0.4,'MDUAL',410.73,17.9,78,41,93),
  @('GAS','GAS_Q100',0.4,'SOP',31853.05,371.1,458,41,93),
  @('GAS','GAS_Q100',0.4,'pMCSKY',31742.56,312.5,396,41,93),
  @('GAS','GAS_Q100',0.8,'MDUAL',516.11,20,80,40,112),
  @('GAS','GAS_Q100',0.8,'SOP',32355.13,363.8,457,40,112),
  @('GAS','GAS_Q100',0.8,'pMCSKY',32128.16,303.9,401,40,112),
  @('GAS','GAS_Q100',1,'MDUAL',454.16,16.7,80,32,86),
  @('GAS','GAS_Q100',1,'SOP',32009.18,361.6,448,32,86),
  @('GAS','GAS_Q100',1,'pMCSKY',31955.34,302.6,392,32,86),
  @('HPC','HPC_Q100',0,'MDUAL',253.8,11.7,31.5,93.5,285.5),
  @('HPC','HPC_Q100',0,'SOP',2211.135,138.75,193.5,93.5,285.5),
  @('HPC','HPC_Q100',0,'pMCSKY',2058.805,96.35,156.5,93.5,285.5),
  @('HPC','HPC_Q100',0.2,'MDUAL',257.924,10.16,31.8,97.2,251.8),
  @('HPC','HPC_Q100',0.2,'SOP',6579.101999999999,142.98,205.8,97.2,251.8),
  @('HPC','HPC_Q100',0.2,'pMCSKY',6093.838,78.42,142.2,97.2,251.8),
  @('HPC','HPC_Q100',0.4,'MDUAL',217.19,10.2,32,91,231),
  @('HPC','HPC_Q100',0.4,'SOP',5830.33,142.85,199,94.5,240.5),
  @('HPC','HPC_Q100',0.4,'pMCSKY',5595.19,77.8,133,91,231),
  @('HPC','HPC_Q100',0.8,'MDUAL',232.54,10.2,33,95,241),
  @('HPC','HPC_Q100',0.8,'SOP',5838.28,143.2,206,95,241),
  @('HPC','HPC_Q100',0.8,'pMCSKY',5603.52,78.1,149,95,241),
  @('HPC','HPC_Q100',1,'MDUAL',230.24,10.1,32,99,243),
  @('HPC','HPC_Q100',1,'SOP',5847.93,143.2,207,99,243),
  @('HPC','HPC_Q100',1,'pMCSKY',5527.89,77.8,142,99,243),
  @('STK','STK_Q100',0,'MDUAL',22.048,4.12,10,97.2,292.8),
  @('STK','STK_Q100',0,'SOP',3286.796,173.74,196.2,97.2,292.8),
  @('STK','STK_Q100',0,'pMCSKY',3156.284,132.36,155.2,97.2,292.8),
  @('STK','STK_Q100',0.2,'MDUAL',23.184,4,11.2,86.2,278.2),
  @('STK','STK_Q100',0.2,'SOP',13180.518,179.58,250.2,86.2,278.2),
  @('STK','STK_Q100',0.2,'pMCSKY',12634.7,121.1,150.8,86.2,278.2),
  @('STK','STK_Q100',0.4,'MDUAL',24.042,4,11.6,88.2,279.4),
  @('STK','STK_Q100',0.4,'SOP',13157.826,179.4,260.6,88.2,279.4),
  @('STK','STK_Q100',0.4,'pMCSKY',12882.892,121.16,156.2,88.2,279.4),
  @('STK','STK_Q100',0.8,'MDUAL',25.112,4,11.4,87.4,277.2),
  @('STK','STK_Q100',0.8,'SOP',12985.364,179.72,207.6,87.4,277.2),
  @('STK','STK_Q100',0.8,'pMCSKY',12801.982,121.64,152,87.4,277.2),
  @('STK','STK_Q100',1,'MDUAL',24.766,4,12,88.6,283),
  @('STK','STK_Q100',1,'SOP',13007.05,179.68,208.2,88.6,283),
  @('STK','STK_Q100',1,'pMCSKY',12784.564,121.66,157,88.6,283),
  @('TAO','TAO_Q100',0,'MDUAL',100.0466666666667,6.099999999999999,16,116,505.6666666666667),
  @('TAO','TAO_Q100',0,'SOP',4565.05,228.1,267,116,505.6666666666667),
  @('TAO','TAO_Q100',0,'pMCSKY',4428.433333333333,186.7666666666667,226,116,505.6666666666667),
  @('TAO','TAO_Q100',0.2,'MDUAL',109.536,6.26,18.4,121.6,560.2),
  @('TAO','TAO_Q100',0.2,'SOP',24179.346,229.84,279.2,121.6,560.2),
  @('TAO','TAO_Q100',0.2,'pMCSKY',23755.652,170.6,220,121.6,560.2),
  @('TAO','TAO_Q100',0.4,'MDUAL',117.3966666666667,6.333333333333333,19,122.3333333333333,564),
  @('TAO','TAO_Q100',0.4,'SOP',24047.51666666667,231.5,330,122.3333333333333,564),
  @('TAO','TAO_Q100',0.4,'pMCSKY',23311.61333333333,172.0666666666667,222.3333333333333,122.3333333333333,564),
  @('TAO','TAO_Q100',0.8,'MDUAL',115.5,6.266666666666667,18,119.3333333333333,554.3333333333334),
  @('TAO','TAO_Q100',0.8,'SOP',23679.91,231.4666666666667,279.6666666666667,119.3333333333333,554.3333333333334),
  @('TAO','TAO_Q100',0.8,'pMCSKY',22931.48333333333,171.4,244.6666666666667,119.3333333333333,554.3333333333334),
  @('TAO','TAO_Q100',1,'MDUAL',121.505,6.2,18.5,119,558.5),
  @('TAO','TAO_Q100',1,'SOP',23665.505,230.7,280.5,119,558.5),
  @('TAO','TAO_Q100',1,'pMCSKY',23988.825,170.7,266.5,119,558.5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}

# Remove the now-obsolete trailing rows (old rows 92-102) entirely
$ws.Range("A92:J102").Clear()
